$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to stay a text value even when the string parses as a
    # number (e.g. "207.32"), matching the original inline-string content.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.991.40"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.560.20"
$ws.Range("E3").Value = "  +0.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
Set-TextValue "D5" "207.32"
$ws.Range("E5").Value = "  +0.24%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.18%  "

# Row 8 - Solana
Set-TextValue "D8" "22.12"
$ws.Range("E8").Value = "  +2.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.08%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.63%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0859"
$ws.Range("E11").Value = "  +0.09%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.782.64"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.543.55"
$ws.Range("E13").Value = "  -0.75%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.02%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.521"
$ws.Range("E15").Value = "  +1.23%  "

# Row 16 - Litecoin
Set-TextValue "D16" "62.05"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.001.27"
$ws.Range("E17").Value = "  +0.42%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "217.17"
$ws.Range("E19").Value = "  +0.81%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +1.98%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.18%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.62%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +0.73%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.87%  "

# Row 25 - Monero
Set-TextValue "D25" "153.18"

# Row 26 - Cosmos
Set-TextValue "D26" "6.64"

# Row 27 - EthereumClassic
Set-TextValue "D27" "15.05"
$ws.Range("E27").Value = "  +1.16%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.46%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.19%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.0468"
$ws.Range("E30").Value = "  +1.07%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.69%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.51%  "

# Rows 33 & 34 - content swap (Maker <-> InternetComputer(DFINITY)); rank (col A) stays put
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.422.42"
$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "3.11"
$ws.Range("E34").Value = "  +3.59%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +2.80%  "

# Row 36 - TrustWalletToken
Set-TextValue "D36" "1.04"
$ws.Range("E36").Value = "  +8.65%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.32"
$ws.Range("E37").Value = "  +1.06%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.59%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "0.530"
$ws.Range("E39").Value = "  +1.28%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +0.34%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.18%  "

# Rows 42, 43, 44 - content rotation (WEMIXToken -> FraxShare -> MXToken -> WEMIXToken); rank (col A) stays put
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.67"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.33"
$ws.Range("E43").Value = "  +2.73%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D44" "0.994"
$ws.Range("E44").Value = "  +0.56%  "

# Row 45 - Aave
Set-TextValue "D45" "64.92"
$ws.Range("E45").Value = "  +1.91%  "

# Row 46 - RenderToken
Set-TextValue "D46" "1.74"
$ws.Range("E46").Value = "  -0.26%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.696.11"
$ws.Range("E47").Value = "  +0.45%  "

# Row 48 - Quant
Set-TextValue "D48" "87.55"
$ws.Range("E48").Value = "  +1.70%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0521"
$ws.Range("E49").Value = "  +0.39%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +0.39%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0958"
$ws.Range("E51").Value = "  -0.21%  "
